# ---------------------------------------------------------------------------
# Nexial json-showcase.xlsx update:
#   - new "aws.ses" command group inserted into the hidden "#system" sheet
#     (new column C => sendMail/sendTextMail; new row A3 => "aws.ses" target)
#   - every existing "#system" column from C..Z shifts right to D..AA
#   - every existing "#system" "target" entry from A3..A26 shifts down to A4..A27
#   - workbook-level defined names updated/added to match
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Write every cell of the sheet to its final (post-edit) value.
$setData = @(
    @{a="A1"; v="target"},
    @{a="B1"; v="aws.s3"},
    @{a="C1"; v="aws.ses"},
    @{a="D1"; v="base"},
    @{a="E1"; v="csv"},
    @{a="F1"; v="desktop"},
    @{a="G1"; v="excel"},
    @{a="H1"; v="external"},
    @{a="I1"; v="image"},
    @{a="J1"; v="io"},
    @{a="K1"; v="jms"},
    @{a="L1"; v="json"},
    @{a="M1"; v="mail"},
    @{a="N1"; v="number"},
    @{a="O1"; v="pdf"},
    @{a="P1"; v="rdbms"},
    @{a="Q1"; v="redis"},
    @{a="R1"; v="sms"},
    @{a="S1"; v="sound"},
    @{a="T1"; v="ssh"},
    @{a="U1"; v="step"},
    @{a="V1"; v="web"},
    @{a="W1"; v="webalert"},
    @{a="X1"; v="webcookie"},
    @{a="Y1"; v="ws"},
    @{a="Z1"; v="ws.async"},
    @{a="AA1"; v="xml"},
    @{a="A2"; v="aws.s3"},
    @{a="B2"; v="assertNotPresent(profile,remotePath)"},
    @{a="C2"; v="sendMail(profile,to,subject,body)"},
    @{a="D2"; v="appendText(var,appendWith)"},
    @{a="E2"; v="compare(expected,actual,failFast)"},
    @{a="F2"; v="assertAttribute(locator,attribute,expected)"},
    @{a="G2"; v="assertPassword(file)"},
    @{a="H2"; v="runJUnit(className)"},
    @{a="I2"; v="compare(baseline,actual)"},
    @{a="J2"; v="assertEqual(expected,actual)"},
    @{a="K2"; v="receive(var,config,waitMs)"},
    @{a="L2"; v="addOrReplace(json,jsonpath,input,var)"},
    @{a="M2"; v="send(profile,to,subject,body)"},
    @{a="N2"; v="assertBetween(num,minNum,maxNum)"},
    @{a="O2"; v="assertContentEqual(actualPdf,expectedPdf)"},
    @{a="P2"; v="resultToCSV(var,csvFile,delim,showHeader)"},
    @{a="Q2"; v="append(profile,key,value)"},
    @{a="R2"; v="sendText(phones,text)"},
    @{a="S2"; v="laser(repeats)"},
    @{a="T2"; v="scpCopyFrom(var,profile,remote,local)"},
    @{a="U2"; v="observe(prompt)"},
    @{a="V2"; v="assertAndClick(locator,label)"},
    @{a="W2"; v="accept()"},
    @{a="X2"; v="assertNotPresent(name)"},
    @{a="Y2"; v="assertReturnCode(var,returnCode)"},
    @{a="Z2"; v="delete(url,body,output)"},
    @{a="AA2"; v="assertCorrectness(xml,schema)"},
    @{a="A3"; v="aws.ses"},
    @{a="B3"; v="assertPresent(profile,remotePath)"},
    @{a="C3"; v="sendTextMail(profile,to,subject,body)"},
    @{a="D3"; v="assertArrayContain(array,expected)"},
    @{a="E3"; v="compareExtended(var,profile,expected,actual)"},
    @{a="F3"; v="assertChecked(name)"},
    @{a="G3"; v="clear(file,worksheet,range)"},
    @{a="H3"; v="runProgram(programPathAndParms)"},
    @{a="I3"; v="convert(source,format,saveTo)"},
    @{a="J3"; v="assertNotEqual(expected,actual)"},
    @{a="K3"; v="sendMap(config,id,payload)"},
    @{a="L3"; v="assertCorrectness(json,schema)"},
    @{a="N3"; v="assertEqual(num1,num2)"},
    @{a="O3"; v="assertFormElementPresent(var,name)"},
    @{a="P3"; v="runFile(var,db,file)"},
    @{a="Q3"; v="assertKeyExists(profile,key)"},
    @{a="S3"; v="play(audio)"},
    @{a="T3"; v="scpCopyTo(var,profile,local,remote)"},
    @{a="U3"; v="perform(instructions)"},
    @{a="V3"; v="assertAttribute(locator,attrName,value)"},
    @{a="W3"; v="assertPresent()"},
    @{a="X3"; v="assertPresent(name)"},
    @{a="Y3"; v="delete(url,body,var)"},
    @{a="Z3"; v="download(url,queryString,saveTo)"},
    @{a="AA3"; v="assertElementCount(xml,xpath,count)"},
    @{a="A4"; v="base"},
    @{a="B4"; v="copyFrom(var,profile,remote,local)"},
    @{a="D4"; v="assertArrayEqual(array1,array2,exactOrder)"},
    @{a="E4"; v="convertExcel(excel,worksheet,csvFile)"},
    @{a="F4"; v="assertDisabled(name)"},
    @{a="G4"; v="clearPassword(file,password)"},
    @{a="I4"; v="crop(image,dimension,saveTo)"},
    @{a="J4"; v="assertReadableFile(file,minByte)"},
    @{a="K4"; v="sendText(config,id,payload)"},
    @{a="L4"; v="assertElementCount(json,jsonpath,count)"},
    @{a="N4"; v="assertGreater(num1,num2)"},
    @{a="O4"; v="assertFormValue(var,name,expected)"},
    @{a="P4"; v="runSQL(var,db,sql)"},
    @{a="Q4"; v="delete(profile,key)"},
    @{a="S4"; v="speak(text)"},
    @{a="T4"; v="sftpCopyFrom(var,profile,remote,local)"},
    @{a="U4"; v="validate(prompt,responses,passResponses)"},
    @{a="V4"; v="assertAttributeContains(locator,attrName,contains)"},
    @{a="W4"; v="assertText(text,matchBy)"},
    @{a="X4"; v="assertValue(name,value)"},
    @{a="Y4"; v="download(url,queryString,saveTo)"},
    @{a="Z4"; v="get(url,queryString,output)"},
    @{a="AA4"; v="assertElementNotPresent(xml,xpath)"},
    @{a="A5"; v="csv"},
    @{a="B5"; v="copyTo(var,profile,local,remote)"},
    @{a="D5"; v="assertArrayNotContain(array,unexpected)"},
    @{a="E5"; v="fromExcel(excel,worksheet,csvFile)"},
    @{a="F5"; v="assertElementPresent(name)"},
    @{a="G5"; v="columnarCsv(file,worksheet,ranges,output)"},
    @{a="I5"; v="resize(image,width,height,saveTo)"},
    @{a="J5"; v="base64(var,file)"},
    @{a="L5"; v="assertElementNotPresent(json,jsonpath)"},
    @{a="N5"; v="assertGreaterOrEqual(num1,num2)"},
    @{a="O5"; v="assertFormValues(var,name,expectedValues,exactOrder)"},
    @{a="P5"; v="runSQLs(var,db,sqls)"},
    @{a="Q5"; v="flushAll(profile)"},
    @{a="S5"; v="speakNoWait(text)"},
    @{a="T5"; v="sftpCopyTo(var,profile,local,remote)"},
    @{a="V5"; v="assertAttributeNotContains(locator,attrName,contains)"},
    @{a="W5"; v="dismiss()"},
    @{a="X5"; v="delete(name)"},
    @{a="Y5"; v="get(url,queryString,var)"},
    @{a="Z5"; v="head(url,output)"},
    @{a="AA5"; v="assertElementPresent(xml,xpath)"},
    @{a="A6"; v="desktop"},
    @{a="B6"; v="delete(var,profile,remotePath)"},
    @{a="D6"; v="assertContains(text,substring)"},
    @{a="F6"; v="assertEnabled(name)"},
    @{a="G6"; v="csv(file,worksheet,range,output)"},
    @{a="J6"; v="compare(expected,actual,failFast)"},
    @{a="L6"; v="assertElementPresent(json,jsonpath)"},
    @{a="N6"; v="assertLesser(num1,num2)"},
    @{a="O6"; v="assertPatternNotPresent(pdf,regex)"},
    @{a="P6"; v="saveResult(db,sql,output)"},
    @{a="Q6"; v="flushDb(profile)"},
    @{a="T6"; v="sftpDelete(var,profile,remote)"},
    @{a="V6"; v="assertAttributeNotPresent(locator,attrName)"},
    @{a="W6"; v="replyCancel(text)"},
    @{a="X6"; v="deleteAll()"},
    @{a="Y6"; v="head(url,var)"},
    @{a="Z6"; v="patch(url,body,output)"},
    @{a="AA6"; v="assertValue(xml,xpath,expected)"},
    @{a="A7"; v="excel"},
    @{a="B7"; v="list(var,profile,remotePath)"},
    @{a="D7"; v="assertCount(text,regex,expects)"},
    @{a="F7"; v="assertHierCells(matchBy,column,expected,nestedOnly)"},
    @{a="G7"; v="json(file,worksheet,range,header,output)"},
    @{a="J7"; v="copyFiles(source,target)"},
    @{a="L7"; v="assertEqual(expected,actual)"},
    @{a="N7"; v="assertLesserOrEqual(num1,num2)"},
    @{a="O7"; v="assertPatternPresent(pdf,regex)"},
    @{a="P7"; v="saveResults(db,sqls,outputDir)"},
    @{a="Q7"; v="rename(profile,current,new)"},
    @{a="T7"; v="sftpList(var,profile,remote)"},
    @{a="V7"; v="assertAttributePresent(locator,attrName)"},
    @{a="W7"; v="replyOK(text)"},
    @{a="X7"; v="save(var,name)"},
    @{a="Y7"; v="header(name,value)"},
    @{a="Z7"; v="post(url,body,output)"},
    @{a="AA7"; v="assertValues(xml,xpath,array,exactOrder)"},
    @{a="A8"; v="external"},
    @{a="B8"; v="moveFrom(var,profile,remote,local)"},
    @{a="D8"; v="assertEmpty(text)"},
    @{a="F8"; v="assertHierRow(matchBy,expected)"},
    @{a="G8"; v="saveData(var,file,worksheet,range)"},
    @{a="J8"; v="count(var,path,pattern)"},
    @{a="L8"; v="assertValue(json,jsonpath,expected)"},
    @{a="N8"; v="average(var,array)"},
    @{a="O8"; v="assertTextArray(pdf,textArray,ordered)"},
    @{a="Q8"; v="set(profile,key,value)"},
    @{a="T8"; v="sftpMoveFrom(var,profile,remote,local)"},
    @{a="V8"; v="assertChecked(locator)"},
    @{a="W8"; v="storeText(var)"},
    @{a="X8"; v="saveAll(var)"},
    @{a="Y8"; v="headerByVar(name,var)"},
    @{a="Z8"; v="put(url,body,output)"},
    @{a="AA8"; v="assertWellformed(xml)"},
    @{a="A9"; v="image"},
    @{a="B9"; v="moveTo(var,profile,local,remote)"},
    @{a="D9"; v="assertEndsWith(text,suffix)"},
    @{a="F9"; v="assertListCount(count)"},
    @{a="G9"; v="saveRange(var,file,worksheet,range)"},
    @{a="J9"; v="deleteFiles(location,recursive)"},
    @{a="L9"; v="assertValues(json,jsonpath,array,exactOrder)"},
    @{a="N9"; v="ceiling(var)"},
    @{a="O9"; v="assertTextNotPresent(pdf,text)"},
    @{a="Q9"; v="store(var,profile,key)"},
    @{a="T9"; v="sftpMoveTo(var,profile,local,remote)"},
    @{a="V9"; v="assertContainCount(locator,text,count)"},
    @{a="Y9"; v="jwtParse(var,token,key)"},
    @{a="AA9"; v="storeCount(xml,xpath,var)"},
    @{a="A10"; v="io"},
    @{a="D10"; v="assertEqual(expected,actual)"},
    @{a="F10"; v="assertLocatorNotPresent(locator)"},
    @{a="G10"; v="setPassword(file,password)"},
    @{a="J10"; v="filter(source,target,matchPattern)"},
    @{a="L10"; v="assertWellformed(json)"},
    @{a="N10"; v="decrement(var,amount)"},
    @{a="O10"; v="assertTextPresent(pdf,text)"},
    @{a="Q10"; v="storeKeys(var,profile,keyPattern)"},
    @{a="V10"; v="assertCssNotPresent(locator,property)"},
    @{a="Y10"; v="jwtSignHS256(var,payload,key)"},
    @{a="AA10"; v="storeValue(xml,xpath,var)"},
    @{a="A11"; v="jms"},
    @{a="D11"; v="assertNotContains(text,substring)"},
    @{a="F11"; v="assertLocatorPresent(locator)"},
    @{a="G11"; v="write(file,worksheet,startCell,data)"},
    @{a="J11"; v="makeDirectory(source)"},
    @{a="L11"; v="fromCsv(csv,header,jsonFile)"},
    @{a="N11"; v="floor(var)"},
    @{a="O11"; v="count(pdf,text,var)"},
    @{a="V11"; v="assertCssPresent(locator,property,value)"},
    @{a="Y11"; v="oauth(var,url,auth)"},
    @{a="AA11"; v="storeValues(xml,xpath,var)"},
    @{a="A12"; v="json"},
    @{a="D12"; v="assertNotEmpty(text)"},
    @{a="F12"; v="assertMenuEnabled(menu)"},
    @{a="G12"; v="writeAcross(file,worksheet,startCell,array)"},
    @{a="J12"; v="moveFiles(source,target)"},
    @{a="L12"; v="storeCount(json,jsonpath,var)"},
    @{a="N12"; v="increment(var,amount)"},
    @{a="O12"; v="saveAsPages(pdf,destination)"},
    @{a="V12"; v="assertElementByAttributes(nameValues)"},
    @{a="Y12"; v="patch(url,body,var)"},
    @{a="A13"; v="mail"},
    @{a="D13"; v="assertNotEqual(expected,actual)"},
    @{a="F13"; v="assertModalDialogNotPresent()"},
    @{a="G13"; v="writeDown(file,worksheet,startCell,array)"},
    @{a="J13"; v="readFile(var,file)"},
    @{a="L13"; v="storeValue(json,jsonpath,var)"},
    @{a="N13"; v="max(var,array)"},
    @{a="O13"; v="saveAsText(pdf,destination)"},
    @{a="V13"; v="assertElementByText(locator,text)"},
    @{a="Y13"; v="post(url,body,var)"},
    @{a="A14"; v="number"},
    @{a="D14"; v="assertStartsWith(text,prefix)"},
    @{a="F14"; v="assertModalDialogPresent()"},
    @{a="G14"; v="writeVar(var,file,worksheet,startCell)"},
    @{a="J14"; v="readProperty(var,file,property)"},
    @{a="L14"; v="storeValues(json,jsonpath,var)"},
    @{a="N14"; v="min(var,array)"},
    @{a="O14"; v="saveFormValues(pdf,var,pageAndLineStartEnd,strategy)"},
    @{a="V14"; v="assertElementCount(locator,count)"},
    @{a="Y14"; v="put(url,body,var)"},
    @{a="A15"; v="pdf"},
    @{a="D15"; v="assertTextOrder(var,descending)"},
    @{a="F15"; v="assertModalDialogTitle(title)"},
    @{a="J15"; v="rename(target,newName)"},
    @{a="N15"; v="round(var,closestDigit)"},
    @{a="O15"; v="saveMetadata(pdf,var)"},
    @{a="V15"; v="assertElementNotPresent(locator)"},
    @{a="Y15"; v="saveResponsePayload(var,file,append)"},
    @{a="A16"; v="rdbms"},
    @{a="D16"; v="assertVarNotPresent(var)"},
    @{a="F16"; v="assertModalDialogTitleByLocator(locator,title)"},
    @{a="J16"; v="saveDiff(var,expected,actual)"},
    @{a="O16"; v="saveToVar(pdf,var)"},
    @{a="V16"; v="assertElementPresent(locator)"},
    @{a="Y16"; v="soap(action,url,payload,var)"},
    @{a="A17"; v="redis"},
    @{a="D17"; v="assertVarPresent(var)"},
    @{a="F17"; v="assertNotChecked(name)"},
    @{a="J17"; v="saveFileMeta(var,file)"},
    @{a="V17"; v="assertFocus(locator)"},
    @{a="Y17"; v="upload(url,body,fileParams,var)"},
    @{a="A18"; v="sms"},
    @{a="D18"; v="clear(vars)"},
    @{a="F18"; v="assertSelected(name,text)"},
    @{a="J18"; v="saveMatches(var,path,filePattern)"},
    @{a="V18"; v="assertFrameCount(count)"},
    @{a="A19"; v="sound"},
    @{a="D19"; v="failImmediate(text)"},
    @{a="F19"; v="assertTableCell(row,column,contains)"},
    @{a="J19"; v="unzip(zipFile,target)"},
    @{a="V19"; v="assertFramePresent(frameName)"},
    @{a="A20"; v="ssh"},
    @{a="D20"; v="incrementChar(var,amount,config)"},
    @{a="F20"; v="assertTableColumnContains(column,contains)"},
    @{a="J20"; v="validate(var,profile,inputFile)"},
    @{a="V20"; v="assertIECompatMode()"},
    @{a="A21"; v="step"},
    @{a="D21"; v="macro(file,sheet,name)"},
    @{a="F21"; v="assertTableContains(contains)"},
    @{a="J21"; v="writeFile(file,content,append)"},
    @{a="V21"; v="assertIENavtiveMode()"},
    @{a="A22"; v="web"},
    @{a="D22"; v="prependText(var,prependWith)"},
    @{a="F22"; v="assertTableRowContains(row,contains)"},
    @{a="J22"; v="writeFileAsIs(file,content,append)"},
    @{a="V22"; v="assertLinkByLabel(label)"},
    @{a="A23"; v="webalert"},
    @{a="D23"; v="repeatUntil(steps,maxWaitMs)"},
    @{a="F23"; v="assertText(name,expected)"},
    @{a="J23"; v="writeProperty(file,property,value)"},
    @{a="V23"; v="assertNotChecked(locator)"},
    @{a="A24"; v="webcookie"},
    @{a="D24"; v="save(var,value)"},
    @{a="F24"; v="assertWindowTitleContains(contains)"},
    @{a="J24"; v="zip(filePattern,zipFile)"},
    @{a="V24"; v="assertNotFocus(locator)"},
    @{a="A25"; v="ws"},
    @{a="D25"; v="saveCount(text,regex,saveVar)"},
    @{a="F25"; v="clear(locator)"},
    @{a="V25"; v="assertNotText(locator,text)"},
    @{a="A26"; v="ws.async"},
    @{a="D26"; v="saveMatches(text,regex,saveVar)"},
    @{a="F26"; v="clearCombo(name)"},
    @{a="V26"; v="assertNotVisible(locator)"},
    @{a="A27"; v="xml"},
    @{a="D27"; v="saveReplace(text,regex,replace,saveVar)"},
    @{a="F27"; v="clearModalDialog(var,button)"},
    @{a="V27"; v="assertOneMatch(locator)"},
    @{a="D28"; v="section(steps)"},
    @{a="F28"; v="clearTextArea(name)"},
    @{a="V28"; v="assertScrollbarHNotPresent(locator)"},
    @{a="D29"; v="split(text,delim,saveVar)"},
    @{a="F29"; v="clearTextBox(name)"},
    @{a="V29"; v="assertScrollbarHPresent(locator)"},
    @{a="D30"; v="startRecording()"},
    @{a="F30"; v="clickButton(name)"},
    @{a="V30"; v="assertScrollbarVNotPresent(locator)"},
    @{a="D31"; v="stopRecording()"},
    @{a="F31"; v="clickByLocator(locator)"},
    @{a="V31"; v="assertScrollbarVPresent(locator)"},
    @{a="D32"; v="substringAfter(text,delim,saveVar)"},
    @{a="F32"; v="clickCheckBox(name)"},
    @{a="V32"; v="assertTable(locator,row,column,text)"},
    @{a="D33"; v="substringBefore(text,delim,saveVar)"},
    @{a="F33"; v="clickExplorerBar(group,item)"},
    @{a="V33"; v="assertText(locator,text)"},
    @{a="D34"; v="substringBetween(text,start,end,saveVar)"},
    @{a="F34"; v="clickFirstMatchRow(nameValues)"},
    @{a="V34"; v="assertTextContains(locator,text)"},
    @{a="D35"; v="verbose(text)"},
    @{a="F35"; v="clickFirstMatchedList(contains)"},
    @{a="V35"; v="assertTextCount(locator,text,count)"},
    @{a="D36"; v="waitFor(waitMs)"},
    @{a="F36"; v="clickIcon(label)"},
    @{a="V36"; v="assertTextList(locator,list,ignoreOrder)"},
    @{a="F37"; v="clickList(row)"},
    @{a="V37"; v="assertTextMatches(text,minMatch,scrollTo)"},
    @{a="F38"; v="clickMenu(menu)"},
    @{a="V38"; v="assertTextNotPresent(text)"},
    @{a="F39"; v="clickOffset(locator,xOffset,yOffset)"},
    @{a="V39"; v="assertTextOrder(locator,descending)"},
    @{a="F40"; v="clickRadio(name)"},
    @{a="V40"; v="assertTextPresent(text)"},
    @{a="F41"; v="clickTab(group,name)"},
    @{a="V41"; v="assertTitle(text)"},
    @{a="F42"; v="clickTableCell(row,column)"},
    @{a="V42"; v="assertValue(locator,value)"},
    @{a="F43"; v="clickTableRow(row)"},
    @{a="V43"; v="assertValueOrder(locator,descending)"},
    @{a="F44"; v="clickTextPane(name,criteria)"},
    @{a="V44"; v="assertVisible(locator)"},
    @{a="F45"; v="clickTextPaneRow(var,index)"},
    @{a="V45"; v="checkAll(locator)"},
    @{a="F46"; v="closeApplication()"},
    @{a="V46"; v="clearLocalStorage()"},
    @{a="F47"; v="collapseHierTable()"},
    @{a="V47"; v="click(locator)"},
    @{a="F48"; v="editCurrentRow(nameValues)"},
    @{a="V48"; v="clickAndWait(locator,waitMs)"},
    @{a="F49"; v="editHierCells(var,matchBy,nameValues)"},
    @{a="V49"; v="clickByLabel(label)"},
    @{a="F50"; v="editTableCells(row,nameValues)"},
    @{a="V50"; v="clickByLabelAndWait(label,waitMs)"},
    @{a="F51"; v="getRowCount(var)"},
    @{a="V51"; v="close()"},
    @{a="F52"; v="hideExplorerBar()"},
    @{a="V52"; v="closeAll()"},
    @{a="F53"; v="login(form,username,password)"},
    @{a="V53"; v="deselect(locator,text)"},
    @{a="F54"; v="maximize()"},
    @{a="V54"; v="deselectMulti(locator,array)"},
    @{a="F55"; v="minimize()"},
    @{a="V55"; v="dismissInvalidCert()"},
    @{a="F56"; v="resize(width,height)"},
    @{a="V56"; v="dismissInvalidCertPopup()"},
    @{a="F57"; v="saveAllTableRows(var)"},
    @{a="V57"; v="doubleClick(locator)"},
    @{a="F58"; v="saveAttributeByLocator(var,locator,attribute)"},
    @{a="V58"; v="doubleClickAndWait(locator,waitMs)"},
    @{a="F59"; v="saveElementCount(var,name)"},
    @{a="V59"; v="doubleClickByLabel(label)"},
    @{a="F60"; v="saveFirstListData(var,contains)"},
    @{a="V60"; v="doubleClickByLabelAndWait(label,waitMs)"},
    @{a="F61"; v="saveFirstMatchedListIndex(var,contains)"},
    @{a="V61"; v="dragAndDrop(fromLocator,toLocator)"},
    @{a="F62"; v="saveHierCells(var,matchBy,column,nestedOnly)"},
    @{a="V62"; v="editLocalStorage(key,value)"},
    @{a="F63"; v="saveHierRow(var,matchBy)"},
    @{a="V63"; v="executeScript(var,script)"},
    @{a="F64"; v="saveListData(var,contains)"},
    @{a="V64"; v="focus(locator)"},
    @{a="F65"; v="saveLocatorCount(var,locator)"},
    @{a="V65"; v="goBack()"},
    @{a="F66"; v="saveModalDialogText(var)"},
    @{a="V66"; v="goBackAndWait()"},
    @{a="F67"; v="saveModalDialogTextByLocator(var,locater)"},
    @{a="V67"; v="maximizeWindow()"},
    @{a="F68"; v="saveProcessId(var,locator)"},
    @{a="V68"; v="mouseOver(locator)"},
    @{a="F69"; v="saveRowCount(var)"},
    @{a="V69"; v="open(url)"},
    @{a="F70"; v="saveTableRows(var,contains)"},
    @{a="V70"; v="openAndWait(url,waitMs)"},
    @{a="F71"; v="saveTableRowsRange(var,beginRow,endRow)"},
    @{a="V71"; v="openHttpBasic(url,username,password)"},
    @{a="F72"; v="saveText(var,name)"},
    @{a="V72"; v="refresh()"},
    @{a="F73"; v="saveTextPane(var,name,criteria)"},
    @{a="V73"; v="refreshAndWait()"},
    @{a="F74"; v="saveWindowTitle(var)"},
    @{a="V74"; v="resizeWindow(width,height)"},
    @{a="F75"; v="scanTable(var,name)"},
    @{a="V75"; v="saveAllWindowIds(var)"},
    @{a="F76"; v="selectCombo(name,text)"},
    @{a="V76"; v="saveAllWindowNames(var)"},
    @{a="F77"; v="sendKeysToTextBox(name,text1,text2,text3,text4)"},
    @{a="V77"; v="saveAttribute(var,locator,attrName)"},
    @{a="F78"; v="showExplorerBar()"},
    @{a="V78"; v="saveCount(var,locator)"},
    @{a="F79"; v="toggleExplorerBar()"},
    @{a="V79"; v="saveDivsAsCsv(headers,rows,cells,nextPage,file)"},
    @{a="F80"; v="typeAppendTextArea(name,text1,text2,text3,text4)"},
    @{a="V80"; v="saveElement(var,locator)"},
    @{a="F81"; v="typeAppendTextBox(name,text1,text2,text3,text4)"},
    @{a="V81"; v="saveElements(var,locator)"},
    @{a="F82"; v="typeByLocator(locator,text)"},
    @{a="V82"; v="saveLocalStorage(var,key)"},
    @{a="F83"; v="typeTextArea(name,text1,text2,text3,text4)"},
    @{a="V83"; v="saveLocation(var)"},
    @{a="F84"; v="typeTextBox(name,text1,text2,text3,text4)"},
    @{a="V84"; v="savePageAs(var,sessionIdName,url)"},
    @{a="F85"; v="useApp(appId)"},
    @{a="V85"; v="savePageAsFile(sessionIdName,url,file)"},
    @{a="F86"; v="useForm(formName)"},
    @{a="V86"; v="saveTableAsCsv(locator,nextPageLocator,file)"},
    @{a="F87"; v="useHierTable(var,name)"},
    @{a="V87"; v="saveText(var,locator)"},
    @{a="F88"; v="useList(var,name)"},
    @{a="V88"; v="saveTextArray(var,locator)"},
    @{a="F89"; v="useTable(var,name)"},
    @{a="V89"; v="saveTextSubstringAfter(var,locator,delim)"},
    @{a="F90"; v="useTableRow(var,row)"},
    @{a="V90"; v="saveTextSubstringBefore(var,locator,delim)"},
    @{a="F91"; v="waitFor(name,maxWaitMs)"},
    @{a="V91"; v="saveTextSubstringBetween(var,locator,start,end)"},
    @{a="F92"; v="waitForLocator(locator,maxWaitMs)"},
    @{a="V92"; v="saveValue(var,locator)"},
    @{a="V93"; v="scrollLeft(locator,pixel)"},
    @{a="V94"; v="scrollRight(locator,pixel)"},
    @{a="V95"; v="scrollTo(locator)"},
    @{a="V96"; v="select(locator,text)"},
    @{a="V97"; v="selectFrame(locator)"},
    @{a="V98"; v="selectMulti(locator,array)"},
    @{a="V99"; v="selectMultiOptions(locator)"},
    @{a="V100"; v="selectText(locator)"},
    @{a="V101"; v="selectWindow(winId)"},
    @{a="V102"; v="selectWindowAndWait(winId,waitMs)"},
    @{a="V103"; v="selectWindowByIndex(index)"},
    @{a="V104"; v="selectWindowByIndexAndWait(index,waitMs)"},
    @{a="V105"; v="toggleSelections(locator)"},
    @{a="V106"; v="type(locator,value)"},
    @{a="V107"; v="typeKeys(locator,value)"},
    @{a="V108"; v="uncheckAll(locator)"},
    @{a="V109"; v="unselectAllText()"},
    @{a="V110"; v="upload(fieldLocator,file)"},
    @{a="V111"; v="verifyContainText(locator,text)"},
    @{a="V112"; v="verifyText(locator,text)"},
    @{a="V113"; v="wait(waitMs)"},
    @{a="V114"; v="waitForElementPresent(locator)"},
    @{a="V115"; v="waitForPopUp(winId,waitMs)"},
    @{a="V116"; v="waitForTextPresent(text)"},
    @{a="V117"; v="waitForTitle(text)"},
)
foreach ($item in $setData) { $ws.Range($item.a).Value = $item.v }

# 2) Blank out cells that are left over from the pre-shift layout and are not
#    part of the post-edit layout (i.e. would otherwise retain stale values).
$clearAddrs = @(
    "M3",
    "R3",
    "C4",
    "H4",
    "M4",
    "R4",
    "C5",
    "H5",
    "K5",
    "M5",
    "R5",
    "U5",
    "C6",
    "E6",
    "I6",
    "K6",
    "M6",
    "S6",
    "U6",
    "C7",
    "E7",
    "I7",
    "K7",
    "M7",
    "S7",
    "U7",
    "C8",
    "E8",
    "I8",
    "K8",
    "M8",
    "P8",
    "S8",
    "U8",
    "C9",
    "E9",
    "I9",
    "K9",
    "M9",
    "P9",
    "S9",
    "U9",
    "X9",
    "Z9",
    "C10",
    "E10",
    "I10",
    "K10",
    "M10",
    "P10",
    "U10",
    "X10",
    "Z10",
    "C11",
    "E11",
    "I11",
    "K11",
    "M11",
    "U11",
    "X11",
    "Z11",
    "C12",
    "E12",
    "I12",
    "K12",
    "M12",
    "U12",
    "X12",
    "C13",
    "E13",
    "I13",
    "K13",
    "M13",
    "U13",
    "X13",
    "C14",
    "E14",
    "I14",
    "K14",
    "M14",
    "U14",
    "X14",
    "C15",
    "E15",
    "I15",
    "M15",
    "U15",
    "X15",
    "C16",
    "E16",
    "I16",
    "N16",
    "U16",
    "X16",
    "C17",
    "E17",
    "I17",
    "U17",
    "X17",
    "C18",
    "E18",
    "I18",
    "U18",
    "C19",
    "E19",
    "I19",
    "U19",
    "C20",
    "E20",
    "I20",
    "U20",
    "C21",
    "E21",
    "I21",
    "U21",
    "C22",
    "E22",
    "I22",
    "U22",
    "C23",
    "E23",
    "I23",
    "U23",
    "C24",
    "E24",
    "I24",
    "U24",
    "C25",
    "E25",
    "U25",
    "C26",
    "E26",
    "U26",
    "C27",
    "E27",
    "U27",
    "C28",
    "E28",
    "U28",
    "C29",
    "E29",
    "U29",
    "C30",
    "E30",
    "U30",
    "C31",
    "E31",
    "U31",
    "C32",
    "E32",
    "U32",
    "C33",
    "E33",
    "U33",
    "C34",
    "E34",
    "U34",
    "C35",
    "E35",
    "U35",
    "C36",
    "E36",
    "U36",
    "E37",
    "U37",
    "E38",
    "U38",
    "E39",
    "U39",
    "E40",
    "U40",
    "E41",
    "U41",
    "E42",
    "U42",
    "E43",
    "U43",
    "E44",
    "U44",
    "E45",
    "U45",
    "E46",
    "U46",
    "E47",
    "U47",
    "E48",
    "U48",
    "E49",
    "U49",
    "E50",
    "U50",
    "E51",
    "U51",
    "E52",
    "U52",
    "E53",
    "U53",
    "E54",
    "U54",
    "E55",
    "U55",
    "E56",
    "U56",
    "E57",
    "U57",
    "E58",
    "U58",
    "E59",
    "U59",
    "E60",
    "U60",
    "E61",
    "U61",
    "E62",
    "U62",
    "E63",
    "U63",
    "E64",
    "U64",
    "E65",
    "U65",
    "E66",
    "U66",
    "E67",
    "U67",
    "E68",
    "U68",
    "E69",
    "U69",
    "E70",
    "U70",
    "E71",
    "U71",
    "E72",
    "U72",
    "E73",
    "U73",
    "E74",
    "U74",
    "E75",
    "U75",
    "E76",
    "U76",
    "E77",
    "U77",
    "E78",
    "U78",
    "E79",
    "U79",
    "E80",
    "U80",
    "E81",
    "U81",
    "E82",
    "U82",
    "E83",
    "U83",
    "E84",
    "U84",
    "E85",
    "U85",
    "E86",
    "U86",
    "E87",
    "U87",
    "E88",
    "U88",
    "E89",
    "U89",
    "E90",
    "U90",
    "E91",
    "U91",
    "E92",
    "U92",
    "U93",
    "U94",
    "U95",
    "U96",
    "U97",
    "U98",
    "U99",
    "U100",
    "U101",
    "U102",
    "U103",
    "U104",
    "U105",
    "U106",
    "U107",
    "U108",
    "U109",
    "U110",
    "U111",
    "U112",
    "U113",
    "U114",
    "U115",
    "U116",
    "U117",
)
foreach ($a in $clearAddrs) { $ws.Range($a).Value = "" }

# 3) Touch AB1 (format-only, no value) so the sheet dimension keeps matching
#    its historical one-column overshoot (A1:AA117 -> A1:AB117).
$ws.Range("AB1").Font.Bold = $ws.Range("AB1").Font.Bold

# 4) Update defined names whose target range shifted, and add the new one.
$nameUpdates = @(
    @{n="base"; r="'#system'!`$D`$2:`$D`$36"},
    @{n="csv"; r="'#system'!`$E`$2:`$E`$5"},
    @{n="desktop"; r="'#system'!`$F`$2:`$F`$92"},
    @{n="excel"; r="'#system'!`$G`$2:`$G`$14"},
    @{n="external"; r="'#system'!`$H`$2:`$H`$3"},
    @{n="image"; r="'#system'!`$I`$2:`$I`$5"},
    @{n="io"; r="'#system'!`$J`$2:`$J`$24"},
    @{n="jms"; r="'#system'!`$K`$2:`$K`$4"},
    @{n="json"; r="'#system'!`$L`$2:`$L`$14"},
    @{n="mail"; r="'#system'!`$M`$2:`$M`$2"},
    @{n="number"; r="'#system'!`$N`$2:`$N`$15"},
    @{n="pdf"; r="'#system'!`$O`$2:`$O`$16"},
    @{n="rdbms"; r="'#system'!`$P`$2:`$P`$7"},
    @{n="redis"; r="'#system'!`$Q`$2:`$Q`$10"},
    @{n="sms"; r="'#system'!`$R`$2:`$R`$2"},
    @{n="sound"; r="'#system'!`$S`$2:`$S`$5"},
    @{n="ssh"; r="'#system'!`$T`$2:`$T`$9"},
    @{n="step"; r="'#system'!`$U`$2:`$U`$4"},
    @{n="target"; r="'#system'!`$A`$2:`$A`$27"},
    @{n="web"; r="'#system'!`$V`$2:`$V`$117"},
    @{n="webalert"; r="'#system'!`$W`$2:`$W`$8"},
    @{n="webcookie"; r="'#system'!`$X`$2:`$X`$8"},
    @{n="ws"; r="'#system'!`$Y`$2:`$Y`$17"},
    @{n="ws.async"; r="'#system'!`$Z`$2:`$Z`$8"},
    @{n="xml"; r="'#system'!`$AA`$2:`$AA`$11"},
)
foreach ($nu in $nameUpdates) { $wb.Names.Item($nu.n).RefersTo = "=" + $nu.r }

$wb.Names.Add("aws.ses", "=" + "'#system'!`$C`$2:`$C`$3")

